$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.018666565760348
$ws.Cells.Item(2, 4).Value = 1.024073826565401
$ws.Cells.Item(2, 5).Value = 0.9926147277508489
$ws.Cells.Item(2, 6).Value = 1.029921335933025
$ws.Cells.Item(2, 9).Value = 1.028572256807639
$ws.Cells.Item(2, 10).Value = 1.023873704226621
$ws.Cells.Item(2, 11).Value = 1.026903397120257
$ws.Cells.Item(2, 12).Value = 0.9955398523336033
$ws.Cells.Item(2, 13).Value = 1.032733845174414
$ws.Cells.Item(2, 14).Value = 1.025327722326579

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.019574927750655
$ws.Cells.Item(3, 4).Value = 1.024712483340146
$ws.Cells.Item(3, 5).Value = 0.9936372048519304
$ws.Cells.Item(3, 6).Value = 1.031052197450709
$ws.Cells.Item(3, 9).Value = 1.028706233435404
$ws.Cells.Item(3, 10).Value = 1.024418702258435
$ws.Cells.Item(3, 11).Value = 1.02734952631437
$ws.Cells.Item(3, 12).Value = 0.9963617723202692
$ws.Cells.Item(3, 13).Value = 1.033672085886387
$ws.Cells.Item(3, 14).Value = 1.025873494318109

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.020162947137686
$ws.Cells.Item(4, 4).Value = 1.025125777582575
$ws.Cells.Item(4, 5).Value = 0.9942998659930995
$ws.Cells.Item(4, 6).Value = 1.031784383025544
$ws.Cells.Item(4, 9).Value = 1.028791592708657
$ws.Cells.Item(4, 10).Value = 1.024771002658053
$ws.Cells.Item(4, 11).Value = 1.027637556739405
$ws.Cells.Item(4, 12).Value = 0.9968940712668345
$ws.Cells.Item(4, 13).Value = 1.034279055542677
$ws.Cells.Item(4, 14).Value = 1.026226295024704

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.020410209189591
$ws.Cells.Item(5, 4).Value = 1.025299534933467
$ws.Cells.Item(5, 5).Value = 0.9945786998346017
$ws.Cells.Item(5, 6).Value = 1.032092299624961
$ws.Cells.Item(5, 9).Value = 1.028827158290757
$ws.Cells.Item(5, 10).Value = 1.024919025266414
$ws.Cells.Item(5, 11).Value = 1.027758489481111
$ws.Cells.Item(5, 12).Value = 0.997117960005301
$ws.Cells.Item(5, 13).Value = 1.034534193024732
$ws.Cells.Item(5, 14).Value = 1.026374527842147

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.020451728995915
$ws.Cells.Item(6, 4).Value = 1.025328710015874
$ws.Cells.Item(6, 5).Value = 0.9946255319796338
$ws.Cells.Item(6, 6).Value = 1.032144006375417
$ws.Cells.Item(6, 9).Value = 1.028833111162685
$ws.Cells.Item(6, 10).Value = 1.024943873938122
$ws.Cells.Item(6, 11).Value = 1.02777878549237
$ws.Cells.Item(6, 12).Value = 0.9971555583673453
$ws.Cells.Item(6, 13).Value = 1.034577029813001
$ws.Cells.Item(6, 14).Value = 1.026399411801819

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.020166250835633
$ws.Cells.Item(7, 4).Value = 1.025128099305195
$ws.Cells.Item(7, 5).Value = 0.9943035907982488
$ws.Cells.Item(7, 6).Value = 1.031788497006759
$ws.Cells.Item(7, 9).Value = 1.028792069193497
$ws.Cells.Item(7, 10).Value = 1.024772980876272
$ws.Cells.Item(7, 11).Value = 1.027639173259693
$ws.Cells.Item(7, 12).Value = 0.9968970624462087
$ws.Cells.Item(7, 13).Value = 1.0342824648305
$ws.Cells.Item(7, 14).Value = 1.02622827605222

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.018973499040638
$ws.Cells.Item(8, 4).Value = 1.024289654443549
$ws.Cells.Item(8, 5).Value = 0.9929600610674301
$ws.Cells.Item(8, 6).Value = 1.030303424144117
$ws.Cells.Item(8, 9).Value = 1.028617810489123
$ws.Cells.Item(8, 10).Value = 1.024057960838395
$ws.Cells.Item(8, 11).Value = 1.027054301618854
$ws.Cells.Item(8, 12).Value = 0.995817528259106
$ws.Cells.Item(8, 13).Value = 1.033050955760219
$ws.Cells.Item(8, 14).Value = 1.025512240603876

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.016873654754303
$ws.Cells.Item(9, 4).Value = 1.022812574271612
$ws.Cells.Item(9, 5).Value = 0.9906006454969559
$ws.Cells.Item(9, 6).Value = 1.027689932325439
$ws.Cells.Item(9, 9).Value = 1.028300557813042
$ws.Cells.Item(9, 10).Value = 1.022795358019926
$ws.Cells.Item(9, 11).Value = 1.026018779895359
$ws.Cells.Item(9, 12).Value = 0.9939188001724441
$ws.Cells.Item(9, 13).Value = 1.030879859065943
$ws.Cells.Item(9, 14).Value = 1.024247844744582

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.015475104874855
$ws.Cells.Item(10, 4).Value = 1.021828175231312
$ws.Cells.Item(10, 5).Value = 0.989033133672735
$ws.Cells.Item(10, 6).Value = 1.025949894216574
$ws.Cells.Item(10, 9).Value = 1.028082232679559
$ws.Cells.Item(10, 10).Value = 1.0219518866512
$ws.Cells.Item(10, 11).Value = 1.025325184151708
$ws.Cells.Item(10, 12).Value = 0.9926553831429383
$ws.Cells.Item(10, 13).Value = 1.02943178917024
$ws.Cells.Item(10, 14).Value = 1.023403175549764

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.014869845405079
$ws.Cells.Item(11, 4).Value = 1.02140201118295
$ws.Cells.Item(11, 5).Value = 0.988355674866747
$ws.Cells.Item(11, 6).Value = 1.025196981959098
$ws.Cells.Item(11, 9).Value = 1.027986082298906
$ws.Cells.Item(11, 10).Value = 1.021586251155603
$ws.Cells.Item(11, 11).Value = 1.025024089164023
$ws.Cells.Item(11, 12).Value = 0.9921088820399291
$ws.Cells.Item(11, 13).Value = 1.028804601612821
$ws.Cells.Item(11, 14).Value = 1.023037020809825

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.014645074079291
$ws.Cells.Item(12, 4).Value = 1.021243729202546
$ws.Cells.Item(12, 5).Value = 0.9881042295826724
$ws.Cells.Item(12, 6).Value = 1.024917397221548
$ws.Cells.Item(12, 9).Value = 1.027950125643365
$ws.Cells.Item(12, 10).Value = 1.021450377387144
$ws.Cells.Item(12, 11).Value = 1.024912135144509
$ws.Cells.Item(12, 12).Value = 0.9919059725120875
$ws.Cells.Item(12, 13).Value = 1.028571611693744
$ws.Cells.Item(12, 14).Value = 1.02290095408503

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.014693286064541
$ws.Cells.Item(13, 4).Value = 1.021277680557051
$ws.Cells.Item(13, 5).Value = 0.9881581567098651
$ws.Cells.Item(13, 6).Value = 1.02497736546002
$ws.Cells.Item(13, 9).Value = 1.02795784941842
$ws.Cells.Item(13, 10).Value = 1.021479525497386
$ws.Cells.Item(13, 11).Value = 1.024936154805252
$ws.Cells.Item(13, 12).Value = 0.9919494934313052
$ws.Cells.Item(13, 13).Value = 1.028621589931469
$ws.Cells.Item(13, 14).Value = 1.022930143588931

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.0148512647341
$ws.Cells.Item(14, 4).Value = 1.021388927245888
$ws.Cells.Item(14, 5).Value = 0.9883348863814464
$ws.Cells.Item(14, 6).Value = 1.025173869757874
$ws.Cells.Item(14, 9).Value = 1.027983115047547
$ws.Cells.Item(14, 10).Value = 1.021575021017205
$ws.Cells.Item(14, 11).Value = 1.025014837329661
$ws.Cells.Item(14, 12).Value = 0.9920921077337197
$ws.Cells.Item(14, 13).Value = 1.028785343093545
$ws.Cells.Item(14, 14).Value = 1.023025774723342

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.014948607117115
$ws.Cells.Item(15, 4).Value = 1.021457471936617
$ws.Cells.Item(15, 5).Value = 0.9884438009545853
$ws.Cells.Item(15, 6).Value = 1.02529495319221
$ws.Cells.Item(15, 9).Value = 1.027998649963413
$ws.Cells.Item(15, 10).Value = 1.021633850966052
$ws.Cells.Item(15, 11).Value = 1.025063301155195
$ws.Cells.Item(15, 12).Value = 0.9921799884222134
$ws.Cells.Item(15, 13).Value = 1.028886233572135
$ws.Cells.Item(15, 14).Value = 1.023084688217464

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.015515280787814
$ws.Cells.Item(16, 4).Value = 1.021856460277478
$ws.Cells.Item(16, 5).Value = 0.9890781214508737
$ws.Cells.Item(16, 6).Value = 1.025999873802948
$ws.Cells.Item(16, 9).Value = 1.028088579888148
$ws.Cells.Item(16, 10).Value = 1.021976144159463
$ws.Cells.Item(16, 11).Value = 1.025345150814696
$ws.Cells.Item(16, 12).Value = 0.9926916645766087
$ws.Cells.Item(16, 13).Value = 1.029473410098565
$ws.Cells.Item(16, 14).Value = 1.023427467506471

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.015870826828263
$ws.Cells.Item(17, 4).Value = 1.022106759337229
$ws.Cells.Item(17, 5).Value = 0.989476357848556
$ws.Cells.Item(17, 6).Value = 1.026442195244144
$ws.Cells.Item(17, 9).Value = 1.028144558670611
$ws.Cells.Item(17, 10).Value = 1.022190747180532
$ws.Cells.Item(17, 11).Value = 1.025521743802507
$ws.Cells.Item(17, 12).Value = 0.9930127773699352
$ws.Cells.Item(17, 13).Value = 1.029841686628819
$ws.Cells.Item(17, 14).Value = 1.023642375288446

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.016078241830166
$ws.Cells.Item(18, 4).Value = 1.022252762879354
$ws.Cells.Item(18, 5).Value = 0.9897087662937556
$ws.Cells.Item(18, 6).Value = 1.026700245552877
$ws.Cells.Item(18, 9).Value = 1.028177054378697
$ws.Cells.Item(18, 10).Value = 1.022315882195121
$ws.Cells.Item(18, 11).Value = 1.025624673791145
$ws.Cells.Item(18, 12).Value = 0.9932001317071769
$ws.Cells.Item(18, 13).Value = 1.030056480295457
$ws.Cells.Item(18, 14).Value = 1.023767688009109

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.016148970226582
$ws.Cells.Item(19, 4).Value = 1.022302547690886
$ws.Cells.Item(19, 5).Value = 0.9897880325774034
$ws.Cells.Item(19, 6).Value = 1.026788242810743
$ws.Cells.Item(19, 9).Value = 1.028188108142593
$ws.Cells.Item(19, 10).Value = 1.022358543315934
$ws.Cells.Item(19, 11).Value = 1.025659757730929
$ws.Cells.Item(19, 12).Value = 0.9932640239640975
$ws.Cells.Item(19, 13).Value = 1.030129716702665
$ws.Cells.Item(19, 14).Value = 1.023810409713607

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.015832676890624
$ws.Cells.Item(20, 4).Value = 1.022079903745947
$ws.Cells.Item(20, 5).Value = 0.9894336180360679
$ws.Cells.Item(20, 6).Value = 1.026394732992989
$ws.Cells.Item(20, 9).Value = 1.02813856878611
$ws.Cells.Item(20, 10).Value = 1.02216772636723
$ws.Cells.Item(20, 11).Value = 1.025502804665138
$ws.Cells.Item(20, 12).Value = 0.9929783193494215
$ws.Cells.Item(20, 13).Value = 1.029802175676428
$ws.Cells.Item(20, 14).Value = 1.023619321782948

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.014804742595943
$ws.Cells.Item(21, 4).Value = 1.021356167463667
$ws.Cells.Item(21, 5).Value = 0.9882828385668249
$ws.Cells.Item(21, 6).Value = 1.025116001921203
$ws.Cells.Item(21, 9).Value = 1.02797568163018
$ws.Cells.Item(21, 10).Value = 1.021546901630663
$ws.Cells.Item(21, 11).Value = 1.024991670426354
$ws.Cells.Item(21, 12).Value = 0.9920501090198102
$ws.Cells.Item(21, 13).Value = 1.028737122545497
$ws.Cells.Item(21, 14).Value = 1.022997615404045

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.01415872326585
$ws.Cells.Item(22, 4).Value = 1.020901209200442
$ws.Cells.Item(22, 5).Value = 0.9875604150241495
$ws.Cells.Item(22, 6).Value = 1.024312478062143
$ws.Cells.Item(22, 9).Value = 1.027871867352052
$ws.Cells.Item(22, 10).Value = 1.021156214967725
$ws.Cells.Item(22, 11).Value = 1.024669641229399
$ws.Cells.Item(22, 12).Value = 0.9914670000341481
$ws.Cells.Item(22, 13).Value = 1.028067339453981
$ws.Cells.Item(22, 14).Value = 1.022606373921234

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.014501163226148
$ws.Cells.Item(23, 4).Value = 1.021142382867521
$ws.Cells.Item(23, 5).Value = 0.9879432794643023
$ws.Cells.Item(23, 6).Value = 1.024738397203886
$ws.Cells.Item(23, 9).Value = 1.027927033902448
$ws.Cells.Item(23, 10).Value = 1.021363358270654
$ws.Cells.Item(23, 11).Value = 1.024840417192564
$ws.Cells.Item(23, 12).Value = 0.991776070289318
$ws.Cells.Item(23, 13).Value = 1.028422417573243
$ws.Cells.Item(23, 14).Value = 1.022813811391413

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.015849915109474
$ws.Cells.Item(24, 4).Value = 1.022092038606022
$ws.Cells.Item(24, 5).Value = 0.9894529299347244
$ws.Cells.Item(24, 6).Value = 1.026416178981395
$ws.Cells.Item(24, 9).Value = 1.028141275838406
$ws.Cells.Item(24, 10).Value = 1.022178128603291
$ws.Cells.Item(24, 11).Value = 1.025511362673904
$ws.Cells.Item(24, 12).Value = 0.9929938892766442
$ws.Cells.Item(24, 13).Value = 1.02982002902378
$ws.Cells.Item(24, 14).Value = 1.023629738791378

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.017416280774737
$ws.Cells.Item(25, 4).Value = 1.023194383745524
$ws.Cells.Item(25, 5).Value = 0.9912096547607049
$ws.Cells.Item(25, 6).Value = 1.028365179487089
$ws.Cells.Item(25, 9).Value = 1.028383779765659
$ws.Cells.Item(25, 10).Value = 1.023122080263795
$ws.Cells.Item(25, 11).Value = 1.026287063108452
$ws.Cells.Item(25, 12).Value = 0.9944092447426414
$ws.Cells.Item(25, 13).Value = 1.03144125909123
$ws.Cells.Item(25, 14).Value = 1.024575030971513

